$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.948866844177246
$ws.Range("B1").Value = 2.534700155258179
$ws.Range("C1").Value = 2.399510383605957
$ws.Range("D1").Value = 2.53086519241333
$ws.Range("E1").Value = 3.286553144454956
